# Applies the "Updated symbol list" price/volume/coin refresh described in the diff.
# Numeric-looking values (columns D, E) are written with a leading apostrophe so the
# interop layer stores them as literal text (matching the original t="inlineStr" cells)
# instead of silently parsing them into numbers/percentages; Style is then reset to
# "Normal" so no visible number-format / quote-prefix styling is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,4).Value = "'307.59"
$ws.Cells.Item(2,4).Style = "Normal"
$ws.Cells.Item(2,5).Value = "'0.79%"
$ws.Cells.Item(2,5).Style = "Normal"
# Row 3
$ws.Cells.Item(3,4).Value = "'36.21"
$ws.Cells.Item(3,4).Style = "Normal"
$ws.Cells.Item(3,5).Value = "'1.24%"
$ws.Cells.Item(3,5).Style = "Normal"
# Row 4
$ws.Cells.Item(4,4).Value = "'5.055"
$ws.Cells.Item(4,4).Style = "Normal"
$ws.Cells.Item(4,5).Value = "'1.01%"
$ws.Cells.Item(4,5).Style = "Normal"
# Row 5
$ws.Cells.Item(5,4).Value = "'0.08109"
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Cells.Item(5,5).Value = "'0.55%"
$ws.Cells.Item(5,5).Style = "Normal"
# Row 6
$ws.Cells.Item(6,4).Value = "'1.995"
$ws.Cells.Item(6,4).Style = "Normal"
$ws.Cells.Item(6,5).Value = "'5.07%"
$ws.Cells.Item(6,5).Style = "Normal"
# Row 7
$ws.Cells.Item(7,2).Value = "KuCoinToken"
$ws.Cells.Item(7,3).Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Cells.Item(7,4).Value = "'7.861"
$ws.Cells.Item(7,4).Style = "Normal"
$ws.Cells.Item(7,5).Value = "'-0.30%"
$ws.Cells.Item(7,5).Style = "Normal"
# Row 8
$ws.Cells.Item(8,2).Value = "MXToken"
$ws.Cells.Item(8,3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(8,4).Value = "'0.9278"
$ws.Cells.Item(8,4).Style = "Normal"
$ws.Cells.Item(8,5).Value = "'-0.25%"
$ws.Cells.Item(8,5).Style = "Normal"
# Row 9
$ws.Cells.Item(9,2).Value = "LiechtensteinCryptoassetsExchange"
$ws.Cells.Item(9,3).Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Cells.Item(9,4).Value = "'0.1457"
$ws.Cells.Item(9,4).Style = "Normal"
$ws.Cells.Item(9,5).Value = "'18.17%"
$ws.Cells.Item(9,5).Style = "Normal"
# Row 10
$ws.Cells.Item(10,2).Value = "WazirX"
$ws.Cells.Item(10,3).Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Cells.Item(10,4).Value = "'0.1930"
$ws.Cells.Item(10,4).Style = "Normal"
$ws.Cells.Item(10,5).Value = "'1.36%"
$ws.Cells.Item(10,5).Style = "Normal"
# Row 11
$ws.Cells.Item(11,2).Value = "MandalaExchangeToken"
$ws.Cells.Item(11,3).Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Cells.Item(11,4).Value = "'0.09136"
$ws.Cells.Item(11,4).Style = "Normal"
$ws.Cells.Item(11,5).Value = "'-0.69%"
$ws.Cells.Item(11,5).Style = "Normal"
# Row 12
$ws.Cells.Item(12,2).Value = "BitrueCoin"
$ws.Cells.Item(12,3).Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Cells.Item(12,4).Value = "'0.03442"
$ws.Cells.Item(12,4).Style = "Normal"
$ws.Cells.Item(12,5).Value = "'-1.99%"
$ws.Cells.Item(12,5).Style = "Normal"
# Row 13
$ws.Cells.Item(13,2).Value = "BitMartToken"
$ws.Cells.Item(13,3).Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Cells.Item(13,4).Value = "'0.09883"
$ws.Cells.Item(13,4).Style = "Normal"
$ws.Cells.Item(13,5).Value = "'-0.38%"
$ws.Cells.Item(13,5).Style = "Normal"
# Row 14
$ws.Cells.Item(14,2).Value = "BitForexToken"
$ws.Cells.Item(14,3).Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Cells.Item(14,4).Value = "'0.001416"
$ws.Cells.Item(14,4).Style = "Normal"
$ws.Cells.Item(14,5).Value = "'-0.34%"
$ws.Cells.Item(14,5).Style = "Normal"
# Row 15
$ws.Cells.Item(15,2).Value = "TigerCash"
$ws.Cells.Item(15,3).Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Cells.Item(15,4).Value = "'0.006334"
$ws.Cells.Item(15,4).Style = "Normal"
$ws.Cells.Item(15,5).Value = "'0.25%"
$ws.Cells.Item(15,5).Style = "Normal"
# Row 16
$ws.Cells.Item(16,2).Value = "LEO"
$ws.Cells.Item(16,3).Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Cells.Item(16,4).Value = "'3.840"
$ws.Cells.Item(16,4).Style = "Normal"
$ws.Cells.Item(16,5).Value = "'6.38%"
$ws.Cells.Item(16,5).Style = "Normal"
# Row 17
$ws.Cells.Item(17,2).Value = "GateToken"
$ws.Cells.Item(17,3).Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Cells.Item(17,4).Value = "'4.165"
$ws.Cells.Item(17,4).Style = "Normal"
$ws.Cells.Item(17,5).Value = "'0.36%"
$ws.Cells.Item(17,5).Style = "Normal"
# Row 18
$ws.Cells.Item(18,4).Value = "'3.451"
$ws.Cells.Item(18,4).Style = "Normal"
$ws.Cells.Item(18,5).Value = "'10.78%"
$ws.Cells.Item(18,5).Style = "Normal"
# Row 19
$ws.Cells.Item(19,4).Value = "'0.3456"
$ws.Cells.Item(19,4).Style = "Normal"
$ws.Cells.Item(19,5).Value = "'0.31%"
$ws.Cells.Item(19,5).Style = "Normal"
# Row 20
$ws.Cells.Item(20,4).Value = "'0.1320"
$ws.Cells.Item(20,4).Style = "Normal"
$ws.Cells.Item(20,5).Value = "'-0.66%"
$ws.Cells.Item(20,5).Style = "Normal"
# Row 21
$ws.Cells.Item(21,4).Value = "'4.824"
$ws.Cells.Item(21,4).Style = "Normal"
$ws.Cells.Item(21,5).Value = "'-6.85%"
$ws.Cells.Item(21,5).Style = "Normal"
# Row 22
$ws.Cells.Item(22,4).Value = "'0.2345"
$ws.Cells.Item(22,4).Style = "Normal"
$ws.Cells.Item(22,5).Value = "'-7.42%"
$ws.Cells.Item(22,5).Style = "Normal"
# Row 23
$ws.Cells.Item(23,4).Value = "'0.04383"
$ws.Cells.Item(23,4).Style = "Normal"
$ws.Cells.Item(23,5).Value = "'-0.84%"
$ws.Cells.Item(23,5).Style = "Normal"
# Row 24
$ws.Cells.Item(24,4).Value = "'0.001234"
$ws.Cells.Item(24,4).Style = "Normal"
$ws.Cells.Item(24,5).Value = "'-0.24%"
$ws.Cells.Item(24,5).Style = "Normal"
# Row 25
$ws.Cells.Item(25,4).Value = "'0.004181"
$ws.Cells.Item(25,4).Style = "Normal"
$ws.Cells.Item(25,5).Value = "'-11.15%"
$ws.Cells.Item(25,5).Style = "Normal"
# Row 27
$ws.Cells.Item(27,4).Value = "'0.0001302"
$ws.Cells.Item(27,4).Style = "Normal"
$ws.Cells.Item(27,5).Value = "'0.12%"
$ws.Cells.Item(27,5).Style = "Normal"
# Row 39
$ws.Cells.Item(39,4).Value = "'0.02048"
$ws.Cells.Item(39,4).Style = "Normal"
$ws.Cells.Item(39,5).Value = "'5.35%"
$ws.Cells.Item(39,5).Style = "Normal"
# Row 40
$ws.Cells.Item(40,4).Value = "'0.05126"
$ws.Cells.Item(40,4).Style = "Normal"
$ws.Cells.Item(40,5).Value = "'-1.12%"
$ws.Cells.Item(40,5).Style = "Normal"
# Row 41
$ws.Cells.Item(41,4).Value = "'0.007468"
$ws.Cells.Item(41,4).Style = "Normal"
$ws.Cells.Item(41,5).Value = "'-0.96%"
$ws.Cells.Item(41,5).Style = "Normal"
# Row 42
$ws.Cells.Item(42,4).Value = "'0.01007"
$ws.Cells.Item(42,4).Style = "Normal"
$ws.Cells.Item(42,5).Value = "'-1.27%"
$ws.Cells.Item(42,5).Style = "Normal"
# Row 43
$ws.Cells.Item(43,4).Value = "'0.1371"
$ws.Cells.Item(43,4).Style = "Normal"
$ws.Cells.Item(43,5).Value = "'0.13%"
$ws.Cells.Item(43,5).Style = "Normal"
# Row 44
$ws.Cells.Item(44,4).Value = "'0.002123"
$ws.Cells.Item(44,4).Style = "Normal"
$ws.Cells.Item(44,5).Value = "'1.07%"
$ws.Cells.Item(44,5).Style = "Normal"
# Row 45
$ws.Cells.Item(45,4).Value = "'0.009864"
$ws.Cells.Item(45,4).Style = "Normal"
$ws.Cells.Item(45,5).Value = "'-8.09%"
$ws.Cells.Item(45,5).Style = "Normal"
# Row 46
$ws.Cells.Item(46,4).Value = "'0.00006306"
$ws.Cells.Item(46,4).Style = "Normal"
$ws.Cells.Item(46,5).Value = "'-0.94%"
$ws.Cells.Item(46,5).Style = "Normal"
# Row 47
$ws.Cells.Item(47,4).Value = "'0.00000000751"
$ws.Cells.Item(47,4).Style = "Normal"
$ws.Cells.Item(47,5).Value = "'-0.09%"
$ws.Cells.Item(47,5).Style = "Normal"
# Row 48
$ws.Cells.Item(48,4).Value = "'63.83"
$ws.Cells.Item(48,4).Style = "Normal"
$ws.Cells.Item(48,5).Value = "'-1.74%"
$ws.Cells.Item(48,5).Style = "Normal"
# Row 49
$ws.Cells.Item(49,4).Value = "'0.001603"
$ws.Cells.Item(49,4).Style = "Normal"
$ws.Cells.Item(49,5).Value = "'-3.61%"
$ws.Cells.Item(49,5).Style = "Normal"
# Row 50
$ws.Cells.Item(50,4).Value = "'0.00002103"
$ws.Cells.Item(50,4).Style = "Normal"
$ws.Cells.Item(50,5).Value = "'-0.09%"
$ws.Cells.Item(50,5).Style = "Normal"
# Row 51
$ws.Cells.Item(51,4).Value = "'0.0002003"
$ws.Cells.Item(51,4).Style = "Normal"
$ws.Cells.Item(51,5).Value = "'-0.09%"
$ws.Cells.Item(51,5).Style = "Normal"
